# Updated cryptos list (prices/volumes) on Sun Sep 29 10:48:45 UTC 2024 with GitHub Actions.
# Rows 28/29 swapped (Kaspa <-> Binance-PegBSC-USD). Numeric-looking price strings are
# forced to remain text (matching the original inline-string cell content) by setting
# NumberFormat to "@" immediately before assigning their Value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.554.51"
$ws.Range("D3").Value = "2.647.11"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.31"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.23"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.127"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.66"
$ws.Range("E13").Value = "  -3.43%  "
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "3.120.35"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "65.376.84"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "2.647.71"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.55"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.19"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.78"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.60"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.69"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.87"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.12"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "527.15"
$ws.Range("E32").Value = "  -3.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.75"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.38"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.419"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.33"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "155.25"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "160.82"
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.05"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0603"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.51"
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.634"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0254"
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0995"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "0.0₆0251"
$ws.Range("E50").Value = "  +6.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.68"
$ws.Range("E51").Value = "  -1.92%  "
